$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 306.5
$ws.Range("I2").Value = 267.8
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 267.8
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -154.8
$ws.Range("N2").Value = -726
$ws.Range("H32").Value = 4526.25
$ws.Range("J32").Value = 4477.4287
$ws.Range("L32").Value = 4477.4287
$ws.Range("N32").Value = -5129.4287
$ws.Range("H76").Value = 8904.538
$ws.Range("J76").Value = 10039.429
$ws.Range("L76").Value = 10039.429
$ws.Range("N76").Value = -10669.429
$ws.Range("H79").Value = 8904.538
$ws.Range("J79").Value = 10039.429
$ws.Range("L79").Value = 10039.429
$ws.Range("N79").Value = -12223.429
$ws.Range("H125").Value = 20976.6
$ws.Range("I125").Value = 32538.666
$ws.Range("J125").Value = 3633.5
$ws.Range("K125").Value = 292847.994
$ws.Range("L125").Value = 32701.5
$ws.Range("M125").Value = -290387.994
$ws.Range("N125").Value = -37621.5
$ws.Range("H132").Value = 8932.933999999999
$ws.Range("J132").Value = 19550
$ws.Range("L132").Value = 58650
$ws.Range("N132").Value = -63710
$ws.Range("H138").Value = 1467.1714
$ws.Range("J138").Value = 1668
$ws.Range("L138").Value = 5004
$ws.Range("N138").Value = -15284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("H124").Value = 26000
$ws.Range("J124").Value = 26000
$ws.Range("L124").Value = 26000
$ws.Range("N124").Value = -35820

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 97
$ws.Range("I22").Value = 99.333336
$ws.Range("J22").Value = 90
$ws.Range("K22").Value = 99.333336
$ws.Range("L22").Value = 90
$ws.Range("M22").Value = 73.666664
$ws.Range("N22").Value = -436
$ws.Range("H95").Value = 90641
$ws.Range("J95").Value = 90641
$ws.Range("L95").Value = 90641
$ws.Range("N95").Value = -96133
$ws.Range("H107").Value = 1939.1578
$ws.Range("I107").Value = 1449.7333
$ws.Range("K107").Value = 1449.7333
$ws.Range("M107").Value = 470.2666999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 563.7778
$ws.Range("I107").Value = 523.0732
$ws.Range("K107").Value = 523.0732
$ws.Range("M107").Value = 1396.9268
$ws.Range("H134").Value = 5439.0386
$ws.Range("I134").Value = 4821.85
$ws.Range("J134").Value = 7496.3335
$ws.Range("K134").Value = 14465.55
$ws.Range("L134").Value = 22489.0005
$ws.Range("M134").Value = -11930.55
$ws.Range("N134").Value = -27559.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 459989.28
$ws.Range("I128").Value = 459989.28
$ws.Range("K128").Value = 1379967.84
$ws.Range("M128").Value = -1374987.84
$ws.Range("H129").Value = 19133536
$ws.Range("I129").Value = 41793144
$ws.Range("J129").Value = 1005848.6
$ws.Range("K129").Value = 125379432
$ws.Range("L129").Value = 3017545.8
$ws.Range("M129").Value = -125374432
$ws.Range("N129").Value = -3027545.8
$ws.Range("H137").Value = 1383.5
$ws.Range("I137").Value = 1188.3334
$ws.Range("K137").Value = 3565.0002
$ws.Range("M137").Value = 1534.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 4386.5
$ws.Range("I31").Value = 4386.5
$ws.Range("K31").Value = 4386.5
$ws.Range("M31").Value = -4094.5
$ws.Range("H37").Value = 4386.5
$ws.Range("I37").Value = 4386.5
$ws.Range("K37").Value = 4386.5
$ws.Range("M37").Value = -4109.5
$ws.Range("H59").Value = 34995.8
$ws.Range("J59").Value = 34995.8
$ws.Range("L59").Value = 34995.8
$ws.Range("N59").Value = -36161.8
$ws.Range("H102").Value = 3358.923
$ws.Range("I102").Value = 2651.4546
$ws.Range("K102").Value = 2651.4546
$ws.Range("M102").Value = -1029.4546

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 41272.637
$ws.Range("I20").Value = 14000
$ws.Range("K20").Value = 14000
$ws.Range("M20").Value = -13774
$ws.Range("H22").Value = 2107.077
$ws.Range("I22").Value = 2461.75
$ws.Range("J22").Value = 1539.6
$ws.Range("K22").Value = 2461.75
$ws.Range("L22").Value = 1539.6
$ws.Range("M22").Value = -2166.75
$ws.Range("N22").Value = -2129.6
$ws.Range("H27").Value = 2107.077
$ws.Range("I27").Value = 2461.75
$ws.Range("J27").Value = 1539.6
$ws.Range("K27").Value = 2461.75
$ws.Range("L27").Value = 1539.6
$ws.Range("M27").Value = -2354.75
$ws.Range("N27").Value = -1753.6
$ws.Range("H42").Value = 16656
$ws.Range("I42").Value = 19999
$ws.Range("K42").Value = 19999
$ws.Range("M42").Value = -19436
$ws.Range("H49").Value = 16656
$ws.Range("I49").Value = 19999
$ws.Range("K49").Value = 19999
$ws.Range("M49").Value = -19852
$ws.Range("H55").Value = 668.3125
$ws.Range("I55").Value = 661
$ws.Range("K55").Value = 661
$ws.Range("M55").Value = -488
$ws.Range("H100").Value = 5464
$ws.Range("I100").Value = 4213.7144
$ws.Range("J100").Value = 6714.2856
$ws.Range("K100").Value = 4213.7144
$ws.Range("L100").Value = 6714.2856
$ws.Range("M100").Value = -3672.7144
$ws.Range("N100").Value = -7796.2856

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 16632.6
$ws.Range("J14").Value = 18570.928
$ws.Range("L14").Value = 18570.928
$ws.Range("N14").Value = -18906.928
$ws.Range("H41").Value = 12162.333
$ws.Range("I41").Value = 11000
$ws.Range("J41").Value = 12394.8
$ws.Range("K41").Value = 11000
$ws.Range("L41").Value = 12394.8
$ws.Range("M41").Value = -10610
$ws.Range("N41").Value = -13174.8
$ws.Range("H96").Value = 1903.25
$ws.Range("I96").Value = 2049.5
$ws.Range("J96").Value = 1854.5
$ws.Range("K96").Value = 2049.5
$ws.Range("L96").Value = 1854.5
$ws.Range("M96").Value = -676.5
$ws.Range("N96").Value = -4600.5
$ws.Range("H105").Value = 42000
$ws.Range("J105").Value = 42000
$ws.Range("L105").Value = 42000
$ws.Range("N105").Value = -48988
$ws.Range("H107").Value = 390.5
$ws.Range("I107").Value = 359.16666
$ws.Range("K107").Value = 1077.49998
$ws.Range("M107").Value = 842.5000199999999

# --- Deletions: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N104").ClearContents()
